$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Delete old row 6 ("Die Fahrzeuge müssen Informationen an den Server senden, wenn sie
# unterwegs sind." / "FA 3"), which shifts the rows below it (old 7,8,9) up to 6,7,8.
$ws.Rows.Item(6).Delete()

# Insert a new blank row so the "Nicht-Funktional" section keeps its original row number.
$ws.Rows.Item(10).Insert()

# Rewrite the requirement texts for the rows that changed wording (in the same order the
# author touched them: GPS text, login text, vehicle-type text, then the two FA ids).
$ws.Range("B6").Value = "Die Fahrzeuge sollen ihre akktuellen GPS Koordinaten an den senden."
$ws.Range("B4").Value = "Die Fahrzeuge sollen dem Server Informationen senden, wenn sie sich anmelden."
$ws.Range("B5").Value = "Die Fahrzeugart muss dem Server mitgeteilt werden."
$ws.Range("A4").Value = "FA2"
$ws.Range("A5").Value = "FA3"

# Update view settings to match
$ws.Range("B10").Select()
